# This workbook tracks daily "Camote" (sweet potato) prices for the
# "Vega Modelo de Temuco" market. The update adds two new daily price
# records (serial dates 44972 and 44973), inserted in chronological
# order within the existing data block. Inserting the rows shifts all
# subsequent rows down by one each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($RowIndex, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($RowIndex, 1).Value2 = 10
    $ws.Cells.Item($RowIndex, 2).Value2 = "Vega Modelo de Temuco"
    $ws.Cells.Item($RowIndex, 3).Value2 = "La Araucanía"
    $ws.Cells.Item($RowIndex, 4).Value2 = $Fecha
    $ws.Cells.Item($RowIndex, 5).Value2 = 9
    $ws.Cells.Item($RowIndex, 6).Value2 = 100114002
    $ws.Cells.Item($RowIndex, 7).Value2 = "Camote"
    $ws.Cells.Item($RowIndex, 8).Value2 = "Sin especificar"
    $ws.Cells.Item($RowIndex, 9).Value2 = "Primera"
    $ws.Cells.Item($RowIndex, 10).Value2 = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value2 = $PrecioMin
    $ws.Cells.Item($RowIndex, 12).Value2 = $PrecioMax
    $ws.Cells.Item($RowIndex, 13).Value2 = $PrecioProm
    $ws.Cells.Item($RowIndex, 14).Value2 = $Unidad
    $ws.Cells.Item($RowIndex, 15).Value2 = $Origen
    $ws.Cells.Item($RowIndex, 16).Value2 = $PrecioKg
    $ws.Cells.Item($RowIndex, 17).Value2 = $KgUnidades
    $ws.Cells.Item($RowIndex, 18).Value2 = "Hortaliza"
}

# Insert first new record at row 44 (pushes the old row 44 and everything
# below it down by one row).
$ws.Rows.Item(44).Insert()
Set-DataRow 44 44972 40 26000 26000 26000 "`$/malla 20 kilos" "Perú" 1300 20

# Insert second new record at row 59 (pushes the row currently at 59 --
# originally row 58 -- and everything below it down by one more row).
$ws.Rows.Item(59).Insert()
Set-DataRow 59 44973 80 26000 26000 26000 "`$/malla 20 kilos" "Perú" 1300 20
